# Facilitator guidelines - Ants Problem: English -> Swahili (Kenya) translation
# Applies the text replacements described by the commit diff, plus updates
# the document's default proofing language from sw-TZ to sw-KE.

$d = $word.ActiveDocument

function Replace-AllText {
    param(
        [string]$Find,
        [string]$Replace
    )
    $d.Content.Find.Execute($Find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $Replace, 2) | Out-Null
}

# Table header / label cells
Replace-AllText "Video Title" "Kichwa cha Video"
Replace-AllText "Topic" "Mada"
Replace-AllText "Aim(s)" "Malengo"
Replace-AllText "Length" "Urefu"
Replace-AllText "Camp Location" "Mahali pa Kambi"
Replace-AllText "Facilitators" "Wawezeshaji"
Replace-AllText "N. of students" "N. ya wanafunzi"
Replace-AllText "Date" "Tarehe"
Replace-AllText "Resources" "Rasilimali"
Replace-AllText "needed" "inahitajika"
Replace-AllText "Preparations" "Maandalizi"
Replace-AllText "Video time" "Muda wa video"
Replace-AllText "What facilitator does" "Mwezeshaji anafanya nini"
Replace-AllText "What learners do" "Wanachofanya wanafunzi"
Replace-AllText "General VMC Video Introduction" "Utangulizi Mkuu wa Video ya VMC"
Replace-AllText "Video Introduction" "Utangulizi wa Video"
Replace-AllText "Riddle" "Kitendawili"
# Appears twice in the document; ReplaceAll (wdReplaceAll = 2) covers both.
Replace-AllText "Assist the process, provoke thoughts" "Kusaidia mchakato, kuchochea mawazo"
Replace-AllText "Solution" "Suluhisho"

# Document default language: Swahili (Tanzania) -> Swahili (Kenya)
$d.Styles("Normal").LanguageID = "sw-KE"

Write-Output "Translation edits applied."
